$d = $word.ActiveDocument

# Locate the "Port  : 8082" run; its paragraph currently also holds the
# trailing _GoBack bookmark.
$r = $d.Content
$r.Find.Execute("Port  : 8082", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$portParaIndex = $r.Paragraphs(1).Index

# Split right after that text so the bookmark ends up alone in its own new
# (trailing) paragraph, separated from "Port  : 8082".
$r.Collapse(0)
$r.InsertBefore([char]13)

# The paragraph that now holds just the bookmark is the next one.
$bookmarkPara = $d.Paragraphs($portParaIndex + 1)
$insertPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)

# Inject the new "Item interface" heading, the TbItem code sample, and the
# trailing blank paragraph as literal OOXML so every run's formatting
# (fonts, bold, colors, shading) matches exactly.
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Item 接口</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="3"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:shd w:val="clear" w:fill="21282D"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E0E2E4"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="6CA3C9"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t xml:space="preserve">TbItem </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E0E2E4"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t xml:space="preserve">item </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t xml:space="preserve">= </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="ECBA61"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>restTemplate</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="F5F6E7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>getForObject</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="FFFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="4CD656"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>"http://MANAGER-SERVICE/item/getById/"</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>+</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E0E2E4"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>id</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="6CA3C9"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>TbItem</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="E784A2"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>class</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="FFFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Consolas" w:hAnsi="Consolas" w:eastAsia="Consolas" w:cs="Consolas"/><w:color w:val="E8E2B7"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:shd w:val="clear" w:fill="21282D"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml)
